$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder header cells D1:G1 (sexo / edad / telefono / correo) ---
$ws.Range("D1").Value = "sexo"
$ws.Range("E1").Value = "edad"
$ws.Range("F1").Value = "telefono"
$ws.Range("G1").Value = "correo"

# --- Highlight the telefono / correo header cells (F1:G1) ---
$hdr = $ws.Range("F1:G1")
$hdr.Style = "Note"
$noteStyle = $wb.Styles.Item("Note")
$noteStyle.Font.Size = 11

# --- Fill in the visit rows (2-6), column by column so shared strings line up
#     the way they were actually authored (whole field filled down first) ---

# A: tipo_documento
$ws.Range("A2").Value = "c"
$ws.Range("A3").Value = "c"
$ws.Range("A4").Value = "c"
$ws.Range("A5").Value = "c"
$ws.Range("A6").Value = "c"

# C: nombre
$ws.Range("C2").Value = "jesus agudo"
$ws.Range("C3").Value = "jesus agudo"
$ws.Range("C4").Value = "jesus agudo"
$ws.Range("C5").Value = "jesus agudo"
$ws.Range("C6").Value = "jesus agudo"

# D: sexo
$ws.Range("D2").Value = "m"
$ws.Range("D3").Value = "m"
$ws.Range("D4").Value = "m"
$ws.Range("D5").Value = "m"
$ws.Range("D6").Value = "m"

# E: edad
$ws.Range("E2").Value = 18
$ws.Range("E3").Value = 19
$ws.Range("E4").Value = 30
$ws.Range("E5").Value = 50
$ws.Range("E6").Value = 23

# F: telefono
$ws.Range("F2").Value = 232323
$ws.Range("F3").Value = 43434
$ws.Range("F4").Value = 34345
$ws.Range("F5").Value = 343443
$ws.Range("F6").Value = 4343432

# G: correo
$ws.Range("G2").Value = "jagudo2514@mgil.aocm"
$ws.Range("G3").Value = "jagudo25@asja.com"
$ws.Range("G4").Value = "sdsdsd@afjjsd.com"
$ws.Range("G5").Value = "samdkasjd@gjjs.com"
$ws.Range("G6").Value = "asmdjasjd@ksdk.com"

# H: provincia
$ws.Range("H2").Value = "Veraguas"
$ws.Range("H3").Value = "Veraguas"
$ws.Range("H4").Value = "Veraguas"
$ws.Range("H5").Value = "Veraguas"
$ws.Range("H6").Value = "Veraguas"

# I: distrito
$ws.Range("I2").Value = "Santiago"
$ws.Range("I3").Value = "Santiago"
$ws.Range("I4").Value = "Santiago"
$ws.Range("I5").Value = "Santiago"
$ws.Range("I6").Value = "Santiago"

# J: corregimiento
$ws.Range("J2").Value = "Santiago"
$ws.Range("J3").Value = "Santiago"
$ws.Range("J4").Value = "Santiago"
$ws.Range("J5").Value = "Santiago"
$ws.Range("J6").Value = "Santiago"

# B: documento (filled last - first row kept as text, rest as plain numbers)
$ws.Range("B2").Value = "9-755-1542"
$ws.Range("B3").Value = 2326
$ws.Range("B4").Value = 2327
$ws.Range("B5").Value = 2358
$ws.Range("B6").Value = 2323

# --- Turn the correo column (emails) into mailto hyperlinks ---
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:jagudo2514@mgil.aocm")
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:jagudo25@asja.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:sdsdsd@afjjsd.com")
$ws.Hyperlinks.Add($ws.Range("G5"), "mailto:samdkasjd@gjjs.com")
$ws.Hyperlinks.Add($ws.Range("G6"), "mailto:asmdjasjd@ksdk.com")

# --- Column widths to fit the new data ---
$ws.Columns.Item(1).ColumnWidth = 16.11
$ws.Columns.Item(2).ColumnWidth = 11.11
$ws.Columns.Item(3).ColumnWidth = 16.66
$ws.Columns.Item(4).ColumnWidth = 5.33
$ws.Columns.Item(5).ColumnWidth = 7.44
$ws.Columns.Item(6).ColumnWidth = 11.78
$ws.Columns.Item(7).ColumnWidth = 33.55
$ws.Columns.Item(8).ColumnWidth = 15.55
$ws.Columns.Item(9).ColumnWidth = 13.78
$ws.Columns.Item(10).ColumnWidth = 14.78

# --- Selection left where the user ended up editing ---
$ws.Range("B3").Select()
